$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Several "Price" values in column D are plain-decimal-looking text
# (e.g. "1.00", "48.40") that must stay exact strings (incl. trailing
# zeros), so NumberFormat is forced to Text ("@") right before those
# writes to stop them being auto-coerced into numbers. Values that already
# contain two dots (thousands separator + decimal, e.g. "43.815.73") are
# never auto-parsed as numbers and don't need that treatment.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.815.73"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.293.43"
$ws.Range("E3").Value = "  -0.01%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.31%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.08"
$ws.Range("E5").Value = "  +15.83%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "269.22"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.13%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.25%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  +1.41%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.40"
$ws.Range("E10").Value = "  +6.43%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0942"
$ws.Range("E11").Value = "  +0.48%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +14.62%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.55%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.79"
$ws.Range("E14").Value = "  -0.23%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.638.40"
$ws.Range("E15").Value = "  +0.04%  "

# Row 16 - Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.308.35"
$ws.Range("E17").Value = "  +0.57%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.714.58"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -1.57%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.10"
$ws.Range("E20").Value = "  +14.28%  "

# Row 21 - Litecoin
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.35"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22 - ImmutableX
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.43"
$ws.Range("E22").Value = "  -0.54%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("E23").Value = "  +8.29%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.86"
$ws.Range("E24").Value = "  -0.24%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +7.26%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.70"
$ws.Range("E27").Value = "  +3.34%  "

# Row 28 - InjectiveProtocol
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "43.73"
$ws.Range("E28").Value = "  +14.55%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  +0.08%  "

# Row 30 - WEMIXToken
$ws.Range("E30").Value = "  -2.18%  "

# Row 31 - Toncoin
$ws.Range("E31").Value = "  -1.35%  "

# Row 32 - Monero
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.17"
$ws.Range("E32").Value = "  -0.82%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0932"
$ws.Range("E33").Value = "  +4.44%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.62"
$ws.Range("E34").Value = "  -1.00%  "

# Row 35 - Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.68"
$ws.Range("E35").Value = "  +4.25%  "

# Row 36 - RenderToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.82"
$ws.Range("E36").Value = "  +2.52%  "

# Row 37 - Stellar
$ws.Range("E37").Value = "  +0.14%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +2.86%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -1.84%  "

# Row 40 - NEARProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  +8.61%  "

# Row 41 - Celestia
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.44"
$ws.Range("E41").Value = "  +18.58%  "

# Row 42 - MultiversX
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.42"
$ws.Range("E42").Value = "  +16.75%  "

# Rows 43 and 44 swap places: Algorand <-> LidoDAOToken,
# with updated price/volume values for each coin.
$ws.Range("B43").Value = "LidoDAOToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  +2.58%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.242"
$ws.Range("E44").Value = "  +2.09%  "

# Row 45 - THORChain
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.35"
$ws.Range("E45").Value = "  +21.70%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.09%  "

# Row 47 - ARBITRUM
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.41"
$ws.Range("E47").Value = "  +2.17%  "

# Row 48 - FraxShare
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.80"
$ws.Range("E48").Value = "  -0.45%  "

# Row 49 - TrustWalletToken
$ws.Range("E49").Value = "  +3.90%  "

# Row 50 - Aave
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.66"
$ws.Range("E50").Value = "  +3.99%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0996"
$ws.Range("E51").Value = "  -2.84%  "
